{"js": "// Change: \"...opportunities that require us to make critical decisions...\"\n// becomes \"...opportunities that require making critical decisions...\"\n// and the Word auto-managed \"_GoBack\" bookmark (marks the most recent edit\n// location) moves from its old spot (end of the \"startup tree\" paragraph)\n// to the new edit point (right after \"making \" in the first paragraph).\n\nconst body = context.document.body;\n\n// 1) Find and replace \"us to make\" -> \"making\" in the first paragraph.\nconst hits = body.search(\"us to make\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find target text 'us to make'\");\n}\n\nconst target = hits.items[0];\ntarget.insertText(\"making\", \"Replace\");\nawait context.sync();\n\n// 2) Move the \"_GoBack\" bookmark: delete it from its old location and add it\n//    right after \"making \" (before \"critical decisions...\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst afterMaking = body.search(\"critical decisions in the face of uncertainty\", { matchCase: true });\nafterMaking.load(\"items\");\nawait context.sync();\n\nif (afterMaking.items.length === 0) {\n  throw new Error(\"Could not find anchor text for bookmark placement\");\n}\n\nconst critRange = afterMaking.items[0].getRange(\"Start\");\ncritRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Change: \"...opportunities that require us to make critical decisions...\"\n# becomes \"...opportunities that require making critical decisions...\"\n# and the auto-managed \"_GoBack\" bookmark (Word's \"last edit location\"\n# marker) moves from its old spot (end of the \"startup tree\" paragraph)\n# to the new edit point, right after \"making \" in the first paragraph.\n\n$d = $word.ActiveDocument\n\n# 1) Replace \"us to make \" with \"making \" in the first paragraph.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"us to make \"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.Execute() | Out-Null\nif (-not $find.Found) {\n    throw \"Could not find target text 'us to make '\"\n}\n$rng.Text = \"making \"\n\n# 2) Move the \"_GoBack\" bookmark to right after the text we just typed.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n$rng.Collapse(0)\n$d.Bookmarks.Add(\"_GoBack\", $rng) | Out-Null\n"}
